$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.383.06"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.18"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6300"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07622"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07738"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.852.64"
$ws.Range("E12").Value = "  -6.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.00001116"
$ws.Range("E13").Value = "  +12.08%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6789"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.66"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "2.103.88"
$ws.Range("E17").Value = "  -7.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.180"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "29.402.99"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.64"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.484"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.30"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1397"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.348"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.463"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.298"
$ws.Range("E30").Value = "  +3.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05591"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.030"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7103"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.582"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "1.239.91"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.416"
$ws.Range("E41").Value = "  +5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9040"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.90"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.86"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000121"
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.158"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4017"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.009"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.681"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  -0.14%  "
